# Initial check-in of translations changes.
#
# The "survey" sheet's E1 header "display.text" becomes "display.prompt.text"
# and the "settings" sheet's C1 header "display.title" becomes
# "display.title.text" (translation-key rename); no other cell content
# changes.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Range("E1").Value = "display.prompt.text"

$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"
